$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7647913694381714
$ws.Range("B1").Value = 1.072787523269653
$ws.Range("C1").Value = 1.464800477027893
$ws.Range("D1").Value = 4.593451023101807
$ws.Range("E1").Value = 2.313761949539185
